$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1 with the same style as the existing header row (copy format from AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for every player row with the team's 1999 season record
$wins = 94
$losses = 68
$ties = 0
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}

Write-Output "Updated AD1:AF50 with season record columns"
